# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps recorded during report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), row 2
$wsOverview.Range("G2").Value = "2017-02-09 09:17:43"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
#              "Correspond Handback DateTime" (L2)
$wsZhCn.Range("H2").Value = "2017-02-09 09:17:19"
$wsZhCn.Range("L2").Value = "2017-02-09 09:18:29"

# de-de sheet: "Correspond Handoff Datetime" (H2) and
#              "Correspond Handback DateTime" (L2)
$wsDeDe.Range("H2").Value = "2017-02-09 09:17:43"
$wsDeDe.Range("L2").Value = "2017-02-09 09:18:57"
